$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (lowercase the header labels)
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

# Add new 4th column: header + value
$ws.Range("D1").Value = "alerttext"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D2").Value = "Customer added successfully"

# Size the new column to fit its content
$ws.Columns.Item(4).EntireColumn.AutoFit()

# Move the active selection, matching the author's final cursor position
$ws.Range("D9").Select()
